$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("urti-elastici")
$ws.Columns("A").Insert()
$ws.Columns("D").Insert()
